# Auto-generated Word COM-interop PowerShell edit script
$d = $word.ActiveDocument

# --- Paragraph 1: title date + paper title (two in-run substring swaps) ---
$d.Content.Find.Execute("13.03.25", $true, $false, $false, $false, $false, $true, 1, $false, "12.03.25", 2) | Out-Null
$d.Content.Find.Execute("SLIM: Let LLM Learn More and Forget Less with Soft LoRA and Identity Mixture", $true, $false, $false, $false, $false, $true, 1, $false, "Transformers are Universal In-context Learner", 2) | Out-Null

# --- Paragraph 2: contains straight double quotes -> set Range.Text directly ---
# (Find/Replace would convert straight quotes to smart quotes, so we avoid it here.)
$d.Paragraphs(2).Range.Text = "היום נסקור קצרות מאמר תיאורטי כבד החוקר את יכולת האקספרסיביות של טרנספורמרים עמוקים. טרנספורמרים הם ארכיטקטורות עמוקות המגדירות ""מיפויים הקשריים"" (in-context mappings), אשר מאפשרים חיזוי של טוקנים חדשים בהתבסס על קבוצת טוקנים נתונה. שימו לב של-in-context כאן יש משמעות קצת שונה מאשר בלמידת in-context הקשור ליכולת של טרנספורמרים ללמוד משימות שלא אומן להם בהתבסס על כמה דוגמאות בפרומפט (לפחות למיטב הבנתי)."

# --- Paragraphs 3-7: full-paragraph text replacement via Find/Replace ---
$d.Content.Find.Execute("המאמר מציע שיטה להלביש ערבוב של מומחים או MoE על LoRa. נזכיר ש-LoRa היא שיטת פיין טיון של רשתות נוירונים שבהם אנו לא מאמנים את כל משקולות המודל אלא רק מטריצות תוספות בעלת ראנק נמוך. MoE היא שיטה להורדה של העומס החישובי בטרנספורמרים כאשר אנו מחלקים את המטריצות בשכבת FFN של הטרנספורמרים לתת-מטריצות (מומחים) כאשר כל פעם לטוקן נתון אנו מפעילים רק חלק מהמומחים. שכבת ניתוב (routing layer) מחשבות את הציון של כל מומחים ובדרך כלל אנו בוחרים k מומחים בעלי ציון הגבוה ביותר (top-k).", $true, $false, $false, $false, $false, $true, 1, $false, "המחברים מוכיחים כי טרנספורמרים עמוקים (בעלי מספר רב של בלוקי הטרנספורמרים) הם מקרבים(approximators) אוניברסליים, כלומר, הם יכולים לקרב כל מיפוי הקשרי רציף מהתפלגויות טוקנים בכל דיוק. יתרה מכך, התוצאות תקפות הן עבור מנגנוני attention דו-כיווניים (כמו באנקודר) והן עבור מנגנוני attention סיבתיים (כמו בדקודרים), תוך שמירה על ממד אמבדינג קבוע שאינו תלוי במספר הטוקנים. ", 2) | Out-Null
$d.Content.Find.Execute("אז המחברים משדכים LoRA עם MoE וזה בדיוק מה שמשך את עיניי. המאמר מציע להחליף LoRA רגיל עם כמה מומחי של LoRA שחלקם הינם מטריצות מראנק 0 או פשוט מטריצות אפסים. לטענת המאמר לא תמיד צריך להפעיל את LoRa. מומחי ה-LoRa נבחרים על ידי רשת ניתוב בדומה ל-MoE הסטנדרתי. עבור כל טוקן נבחרים K מומחים (בינם גם מומחי זהות) בעלי ציונים הגבוהים ביותר. שימו לב שבמאמר יש כמה שגיאות בנוסחאות המחשבים את התוצאה של המנגנון המוצע. ", $true, $false, $false, $false, $false, $true, 1, $false, "הגישה המוצע מבוססת על תורת המידה(סוף סוף מצאתי לה שימוש במאמרי DL), שבה רצפי סדרות מיוצגים כהתפלגויות הסתברותיות במרחב האמבדינגס. זה מאפשר שימוש בכלים מאנליזה פונקציונלית(פלאשבקים מלפני כמעט 30 שנה בתואר הראשון) ובתורת הטרנספורט האופטימלי (כתבתי על זה לא מעט בזמו בהקשר של Wasserstein GAN) על מנת להוכיח את יכולת הקירוב האוניברסלית של טרנספורמרים. תרומה טכנית מרכזית היא הגדרה מחדש של מנגנון ה-attention בטרנספורמרים כאופרטור על התפלגויות. זה מאפשר שימוש במשפט סטון-ויירשטראס(Stone–Weierstrass המהווה הכללה קשוחה של משפט Weierstrass הנלמד באינפי2 לדעתי) - תוצאה יסודית בתורת הקירוב על כך שניתן לקרב כל פונקציה ״נוחה״ על ידי משפחת פונקציות צפופות יחסית (המשפט באמת קשוח המגדיר פונקציות במרחבי האוסדורף וכאלו).", 2) | Out-Null
$d.Content.Find.Execute("לאחר מכן המאמר מציע שיטה לשכלול הציונים של שכבת הניתוב בהתבסס על הסטטיסטיקות של הדאטהסט עליו בוצע הפיינטיון עם השיטה. סטטיסטיקה במקרה הזה מחושבת על המצבים החבויים של הרשת המחושבים על הדאטה של הפיין טיון (אופן החישוב המדויק לא מוגדר בצורה ברורה ולדעתי יש שגיאות בנוסחאות המגדירות אותו). המחברים מציעים לקלסטר את המצבים החבויים האלו לקלסטרים שמספרם כנראה שווה למספר הטוקנים בפרומפט (מוגדר כקבוע במאמר ועבור סדרות קצרות יותר משתמשים בטוקני ה-padding). ", $true, $false, $false, $false, $false, $true, 1, $false, " ייצוג מבוסס-מידה של למידה בהקשר", 2) | Out-Null
$d.Content.Find.Execute("מרכזי הקלסטרים מתעדכנים במהלך הפיין טיון (כל פרומפט הקלט משויך לקלסטר הקרוב ביותר ואז מרכז הקלסטר מחושב מחדש). במהלך האינפרנס פרומפט הקלט משויך לקלסטר הקרוב ביותר (מרחק ריבוע) ואז ציוני המומחים המופקים על ידי שכבת הניתוב עבור מומחי הזהות מוזזים במקדם שעולה אם המרחק לקלסטר הקרוב עולה כאשר הציונים למומחי LoRA האחרים נותרים ללא שינוי. נציין שמרכזי הקלסטרים לא מתעדכנים במהלך האינפרנס.", $true, $false, $false, $false, $false, $true, 1, $false, "חידוש מרכזי במאמר הוא ייצוג של מנגנון ה-attention כאופרטור על התפלגויות במקום על סדרות טוקנים סופיות. דבר זה מאפשר ניתוח אחיד של למידת ההקשר (in-context learning), ללא תלות במספר הטוקנים בסדרה. במקום לעבוד עם קבוצות סופיות של האמבדינגס של הטוקנים, המחברים מגדירים מרחב של התפלגויות הסתברותיות על תת-קבוצה קומפקטית של מרחב אוקלידי (של האמבדינגס). התפלגות משייכת משקלים לאמבדינגס שונים של טוקנים, ובכך מייצגת את המעבר מלמידה על מספר טוקנים סופי לייצוג רציף ואינסופי.", 2) | Out-Null
$d.Content.Find.Execute("לבסוף המאמר מציע דרך לשלב כמה MoE עם LoRa עבור כמה משימות פיין טיון שונות אבל אחרי שגיליתי טעיות גם בפרק הזה, ויתרתי….", $true, $false, $false, $false, $false, $true, 1, $false, "באופן פורמלי, רצף של טוקנים ניתן לייצוג כהתפלגות הסתברות בדידה, המורכבת מסכום משוקלל של פונקציות דלתא דיראק, שכל אחת מהן ממוקמת על הטמעה של טוקן בודד. כאשר מספר הטוקנים גדל, התפלגויות אלה מתכנסות להתפלגויות רציפות. ניסוח זה מאפשר להוכיח תוצאות החלות על כל מספר אפשרי של טוקנים, כולל אינסוף.", 2) | Out-Null

# --- Paragraph 8: URL paragraph text replaced, then 16 new paragraphs appended ---
$d.Content.Find.Execute("https://arxiv.org/pdf/2410.07739", $true, $false, $false, $false, $false, $true, 1, $false, "הגדרת attention כאופרטור על מרחב מידות", 2) | Out-Null

# Append the new paragraphs after the (former URL, now index 8) paragraph, preserving order
$insertAfterIndex = 8
$anchor = $d.Paragraphs($insertAfterIndex).Range
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(9).Range
$anchor.Text = "שכבת טרנספורמר טיפוסית מורכבת משני רכיבים:"
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(10).Range
$anchor.Text = "1. מנגנון attention רב-ראשי, האחראי על עדכון הייצוגים של הטוקנים על ידי חישוב יחסי הגומלין ביניהם."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(11).Range
$anchor.Text = "2. שכבות FFN, המעבדות כל טוקן באופן עצמאי לאחר שלב ה-attention."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(12).Range
$anchor.Text = "המחברים מנסחים מחדש את מנגנון ה-attention כמיפוי הפועל על התפלגויות של טוקנים. במקום לחבר סכום על קבוצת טוקנים סופית, ה-attention מוגדרת כאופרטור אינטגרלי על מרחב ההתפלגויות, מה שהופך את הטוקנים למבנה רציף. ניסוח זה חשוב במיוחד, מכיוון שהוא מאפשר להגדיר רציפות וחלקות של מיפויים בהקשר באמצעות מרחק וסרשטיין (מקרה פרטי שלו הוא earth mover distance), המודד את המרחק בין התפלגויות הסתברותיות. פונקציה היא רציפה במובן וסרשטיין אם שינויים קטנים בהתפלגות הקלט מובילים לשינויים קטנים בהתפלגות הפלט. תכונה זו מבטיחה שהמיפויים שיוצרים טרנספורמרים יציבים לשינויים בהקשר הלימודי."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(13).Range
$anchor.Text = "הוכחת אוניברסליות: קירוב מיפויים הקשריים "
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(14).Range
$anchor.Text = "התוצאות המרכזיות של המאמר מוכיחות כי טרנספורמרים הם מקרבים אוניברסליים למיפויים הקשריים. המחברים מראים כי עבור כל פונקציה רציפה הממפה התפלגויות טוקנים לפלטים, קיים טרנספורמר עמוק שיכול לקרב אותה בכל דיוק. חלק מרכזי בהוכחה הוא בנייה של פונקציות יסודיות בהקשר, המשמשות כיחידות הבסיס לקירוב כל פונקציה כללית במרחבים שהגדרנו קודם. פונקציות אלו הן גרסאות פשוטות יותר של שכבות טרנספורמר, אשר לוכדות את העקרונות המרכזיים של מנגנוני ה-attention. "
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(15).Range
$anchor.Text = "פונקציה יסודית כזו מורכבת משלושה מרכיבים:"
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(16).Range
$anchor.Text = "1. טרנספורמציה לינארית על הטמעת הטוקן (מיפוי אפיני)."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(17).Range
$anchor.Text = "2. אינטראקציה לא-ליניארית המתחשבת בהתפלגות של כלל הטוקנים."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(18).Range
$anchor.Text = "3. התאמה תלוית-הקשר, המאפשרת למודל ""ללמוד בהקשר""."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(19).Range
$anchor.Text = "פונקציות אלו פועלות באופן דומה למנגנון ה-attention  בעל ראש בודד, אך הן קלות יותר לניתוח מתמטי. המחברים מוכיחים כי על ידי הרכבת מספר שכבות של פונקציות אלו, ניתן ליצור טרנספורמרים עמוקים המסוגלים לקרב כל פונקציה בהקשר."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(20).Range
$anchor.Text = "שימוש במשפט סטון-ויירשטראס"
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(21).Range
$anchor.Text = "כדי להוכיח אוניברסליות, המחברים מראים כי קבוצת הפונקציות היסודיות שהם הגדירו מקיימת את תנאי משפט סטון-ויירשטראס, שכאמור הוא משפט מרכזי באנליזה פונקציונלית. המחברים מוכיחים כי הפונקציות היסודיות שלהם מקיימות תנאים אלו, מה שמבטיח כי טרנספורמרים עמוקים יכולים לקרב כל מיפוי הקשרי."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(22).Range
$anchor.Text = "סיכום:"
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(23).Range
$anchor.Text = "המאמר מספק מסגרת מתמטית להוכחת האקספרסיביות של טרנספורמרים בלמידת  מיפוים הקשריים, תוך שימוש באנליזה פונקציונלית, תורת המידה ותורת הטרנספורט האופטימלי. התוצאות מראות כי טרנספורמרים עמוקים יכולים לקרב כל פונקציה תלויה-הקשר, ללא תלות במספר הטוקנים בחלון ההקשר."
$anchor.InsertParagraphAfter()
$anchor = $d.Paragraphs(24).Range
$anchor.Text = "https://arxiv.org/abs/2408.01367"

Write-Output ("FinalParagraphCount=" + $d.Paragraphs.Count)
